$wb = $excel.ActiveWorkbook

# The "ValidLoginCredentials" sheet is the active sheet (A1 currently holds
# the shared string "admin"); update it to "Admin" and move the selection
# from C8 to A8, matching the author's edit.
$ws = $wb.Worksheets.Item("ValidLoginCredentials")
$ws.Activate()

$ws.Range("A1").Value = "Admin"

$ws.Range("A8").Select()
